$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9, shifting existing rows 9-17 down to 10-18
$ws.Rows.Item(9).Insert()

# Update VinMAx value from 35 to 32
$ws.Range("B2").Value = 32

# Fill in the new row 9 ("L ripple current")
$ws.Range("A9").Value = "L ripple current"
$ws.Range("B9").Formula = "=B3*(B2-B3)/(B2*B8*B11)"
$ws.Range("C9").Value = "Amps"

# Apply number format (numFmtId 2 -> "0.00") to B9
$ws.Range("B9").NumberFormat = "0.00"

# Set the active selection to B9
$ws.Range("B9").Select()
